$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    # Force a purely-numeric-looking string to be stored as text rather
    # than being auto-coerced into a number by Excel's input parser.
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

# New column F: header "Zrealizowana" (copy header formatting from A1)
$ws.Range("F1").Value = "Zrealizowana"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# Update row 2 - split purchase into "na recepte" (prescription) vs "bez recepty" (OTC)
Set-TextValue $ws.Range("A2") "202"
$ws.Range("B2").Value = "Lisinopril"
$ws.Range("C2").Value = "BLQY1A35"
$ws.Range("D2").Value = "2025-06-09 11:25:08"
$ws.Range("E2").Value = "2025-07-09 11:25:08"
$ws.Range("F2").Value = "NIE"

$wb.Save()
